# edit.ps1 - applies the "first draft needs changes" revision to
# "LSDE Report Part 1" per the supplied OOXML diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...I think that all of the personal data..." -> "...sensitive data..."
#    (word swapped; the surrounding run is split into three runs in the
#    target, which we approximate by toggling a character property on
#    the replacement text so the engine keeps it as its own run)
# ---------------------------------------------------------------------
$rngPersonal = $d.Content
$foundPersonal = $rngPersonal.Find.Execute("personal ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPersonal) {
    $rngPersonal.Text = "sensitive "
    $rngPersonal.Bold = 1
    $rngPersonal.Bold = 0
}

# ---------------------------------------------------------------------
# 2) Append a new sentence to the end of the "...OpEx model" paragraph.
# ---------------------------------------------------------------------
$opexPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -match "primarily uses the\s*$" -or $candidate.Range.Text -match "OpEx model\s*$") {
        $opexPara = $candidate
    }
}
$opexPara.Range.InsertAfter(" meaning that some financial institutions are weary of using it. ")
$opexIndex = $opexPara.Index

# ---------------------------------------------------------------------
# 3) Insert 3 blank paragraphs, a new "Finally, ..." paragraph, a blank
#    paragraph, and a new "*conclusion..." paragraph right after it.
# ---------------------------------------------------------------------
$cur = $d.Paragraphs.Item($opexIndex)
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($opexIndex + 1)
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Item($opexIndex + 2)
$cur.Range.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($opexIndex + 3)
$cur.Range.InsertParagraphAfter()
$finallyPara = $d.Paragraphs.Item($opexIndex + 4)
$finallyPara.Range.Text = "Finally, a scenario where a financial institution should not use the public cloud for data processing may not be a technological one, it may be because they risk changing their entire landscape and infrastructure that they have built up for multiple years. When switching to a public cloud you are changing the way the interact with your customers and how your employees operate on their day to day. "

$cur = $d.Paragraphs.Item($opexIndex + 4)
$cur.Range.InsertParagraphAfter()

$cur = $d.Paragraphs.Item($opexIndex + 5)
$cur.Range.InsertParagraphAfter()
$conclusionPara = $d.Paragraphs.Item($opexIndex + 6)
$conclusionPara.Range.Text = "*conclusion -hybrid cloud structure* "

# ---------------------------------------------------------------------
# 4) At the very end of the document, add a blank paragraph and a new
#    paragraph containing a hyperlink to the QualiTest article.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()
$newHyperlinkPara = $d.Paragraphs.Last

$qtUrl = "https://qualitestgroup.com/insights/blog/cloud-migration-3-biggest-risks-banks-and-financial-services-companies-need-to-know-and-how-to-avoid-them/"
$d.Hyperlinks.Add($newHyperlinkPara.Range, $qtUrl, [Type]::Missing, [Type]::Missing, $qtUrl)

$trailPara = $d.Paragraphs.Last
$trailPara.Range.InsertAfter(" ")

Write-Output "done"
